# Updated cryptos list on Tue Dec 26 16:27:11 UTC 2023 with GitHub Actions
#
# The coin-ranking list refreshed: rows 2-28 keep the same coin in place
# and only get new Price (D) / Volume(1h) (E) readings, while rows 29-51
# shift down by one slot (a new "LEO" entry was inserted at row 29,
# pushing every coin from "InjectiveProtocol" onward down one row and
# dropping the previous last row, "Cronos") - each row also picks up a
# freshly refreshed Price/Volume reading, not just a copy of the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so numeric-looking
# strings (e.g. "42.459.49", "1.01") are kept as text, matching the
# original inline-string cell type rather than being parsed as numbers.
function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "42.459.49"
Set-TextValue "E2" "  -3.11%  "

Set-TextValue "D3" "2.222.45"
Set-TextValue "E3" "  -2.65%  "

Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.46%  "

Set-TextValue "D5" "111.44"
Set-TextValue "E5" "  -9.96%  "

Set-TextValue "D6" "296.83"
Set-TextValue "E6" "  +11.49%  "

Set-TextValue "D7" "0.625"
Set-TextValue "E7" "  -2.29%  "

Set-TextValue "E8" "  +0.13%  "

Set-TextValue "D9" "0.613"
Set-TextValue "E9" "  -2.17%  "

Set-TextValue "D10" "45.44"
Set-TextValue "E10" "  -6.96%  "

Set-TextValue "D11" "0.0925"
Set-TextValue "E11" "  -1.54%  "

Set-TextValue "D12" "55.20"
Set-TextValue "E12" "  +1.19%  "

Set-TextValue "D13" "8.92"
Set-TextValue "E13" "  -2.91%  "

Set-TextValue "E14" "  -3.40%  "

Set-TextValue "D15" "0.933"
Set-TextValue "E15" "  +4.01%  "

Set-TextValue "D16" "15.17"
Set-TextValue "E16" "  -2.21%  "

Set-TextValue "D17" "2.563.81"
Set-TextValue "E17" "  -2.50%  "

Set-TextValue "D18" "2.253.44"
Set-TextValue "E18" "  -1.07%  "

Set-TextValue "D19" "42.413.45"
Set-TextValue "E19" "  -3.03%  "

Set-TextValue "D20" "7.34"
Set-TextValue "E20" "  +4.69%  "

Set-TextValue "E21" "  -2.65%  "

Set-TextValue "D22" "73.32"
Set-TextValue "E22" "  +1.12%  "

Set-TextValue "D23" "3.56"
Set-TextValue "E23" "  +23.79%  "

Set-TextValue "E24" "  -6.38%  "

Set-TextValue "D25" "230.21"
Set-TextValue "E25" "  -2.49%  "

Set-TextValue "D26" "9.38"
Set-TextValue "E26" "  -2.37%  "

Set-TextValue "D27" "11.81"
Set-TextValue "E27" "  -0.08%  "

Set-TextValue "E28" "  -1.43%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D29" "3.91"
Set-TextValue "E29" "  -1.37%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "38.94"
Set-TextValue "E30" "  -8.75%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D31" "2.22"
Set-TextValue "E31" "  -1.44%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D32" "3.25"
Set-TextValue "E32" "  -3.70%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D33" "173.91"
Set-TextValue "E33" "  +0.57%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "21.20"
Set-TextValue "E34" "  -2.42%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0889"
Set-TextValue "E35" "  -2.60%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "5.69"
Set-TextValue "E36" "  -1.34%  "

Set-TextValue "D37" "4.33"
Set-TextValue "E37" "  +4.52%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "4.82"
Set-TextValue "E38" "  +2.40%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.127"
Set-TextValue "E39" "  -2.19%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.0368"
Set-TextValue "E40" "  -3.14%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.103"
Set-TextValue "E41" "  -3.68%  "

$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D42" "2.50"
Set-TextValue "E42" "  -1.74%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D43" "0.238"
Set-TextValue "E43" "  -0.73%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D44" "71.54"
Set-TextValue "E44" "  -5.97%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D45" "13.17"
Set-TextValue "E45" "  -6.41%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D46" "1.00"
Set-TextValue "E46" "  +0.27%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D47" "1.32"
Set-TextValue "E47" "  -4.15%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D48" "5.50"
Set-TextValue "E48" "  -6.14%  "

$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D49" "1.31"
Set-TextValue "E49" "  +3.27%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D50" "104.46"
Set-TextValue "E50" "  +2.34%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D51" "8.56"
Set-TextValue "E51" "  -2.21%  "
